$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wrap the two JSON condition strings in B2/B3 with [ ... ] to turn them into JSON arrays.
$ws.Range("B2").Value = "[{`n`t""message"": ""[LDAP: error code 49 - 80090308: LdapErr: DSID-0C09042F, comment: AcceptSecurityContext error, data 531, v2580 ]""`n}]"
$ws.Range("B3").Value = "[{`n`t""message"": ""Security token is invalid. java.util.NoSuchElementException: No value present""`n}]"

# Remove the trailing empty row (row 4) that only carried a custom height.
$ws.Rows.Item(4).Delete()

# Move the active selection to B2 (was A3).
$ws.Range("B2").Select()
